$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Overview sheet: the "Ready for handoff" status text becomes
#    "Handed back: in sync with en-US" for both locale columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------------
# 2) zh-cn sheet: handback finished -> fill in Latest Target File / Latest
#    Handback File columns (I, J) and update the Status column (C).
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"

$wsZh.Range("I2").Value = "b69c8e7a-f865-468d-a787-d2e94bd29e7b.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2c290ca3c7e2e838b45c6f337c13a32aa25b2d0/e2e/b69c8e7a-f865-468d-a787-d2e94bd29e7b.md", "", "", "b69c8e7a-f865-468d-a787-d2e94bd29e7b.md")
$wsZh.Range("J2").Value = "b69c8e7a-f865-468d-a787-d2e94bd29e7b.aec7889e3a54a78bddaadf985cbc51a57057d282.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-26 04:28:42"

$wsZh.Range("I3").Value = "fab483d1-3443-49ea-9d8b-da8109a2fa6c.md"
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2c290ca3c7e2e838b45c6f337c13a32aa25b2d0/e2e/fab483d1-3443-49ea-9d8b-da8109a2fa6c.md", "", "", "fab483d1-3443-49ea-9d8b-da8109a2fa6c.md")
$wsZh.Range("J3").Value = "fab483d1-3443-49ea-9d8b-da8109a2fa6c.1bc58dc2306fcc8ec40e6c1211538441ad54691a.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-26 04:28:42"

$wsZh.Columns("C").ColumnWidth = 29.166666666666664
$wsZh.Columns("I").ColumnWidth = 39.16666666666667
$wsZh.Columns("J").ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# 3) de-de sheet: same handback completion.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDe.Range("I2").Value = "b69c8e7a-f865-468d-a787-d2e94bd29e7b.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2c290ca3c7e2e838b45c6f337c13a32aa25b2d0/e2e/b69c8e7a-f865-468d-a787-d2e94bd29e7b.md", "", "", "b69c8e7a-f865-468d-a787-d2e94bd29e7b.md")
$wsDe.Range("J2").Value = "b69c8e7a-f865-468d-a787-d2e94bd29e7b.aec7889e3a54a78bddaadf985cbc51a57057d282.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-26 04:28:48"

$wsDe.Range("I3").Value = "fab483d1-3443-49ea-9d8b-da8109a2fa6c.md"
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/e2c290ca3c7e2e838b45c6f337c13a32aa25b2d0/e2e/fab483d1-3443-49ea-9d8b-da8109a2fa6c.md", "", "", "fab483d1-3443-49ea-9d8b-da8109a2fa6c.md")
$wsDe.Range("J3").Value = "fab483d1-3443-49ea-9d8b-da8109a2fa6c.1bc58dc2306fcc8ec40e6c1211538441ad54691a.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-26 04:28:48"

$wsDe.Columns("C").ColumnWidth = 29.166666666666664
$wsDe.Columns("I").ColumnWidth = 39.16666666666667
$wsDe.Columns("J").ColumnWidth = 39.16666666666667

# ---------------------------------------------------------------------------
# 4) Overview Status columns (E = zh-cn, F = de-de) need to widen to fit the
#    longer status text.
# ---------------------------------------------------------------------------
$wsOverview.Columns("E").ColumnWidth = 29.166666666666664
$wsOverview.Columns("F").ColumnWidth = 29.166666666666664
